$wb = $excel.ActiveWorkbook

# Overview sheet: rows 8 and 9 move from "Ready for handoff" to "In Translation"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

# zh-cn sheet: Status column (C) rows 8 and 9 move to "In Translation"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("C9").Value = "In Translation"

# de-de sheet: Status column (C) rows 8 and 9 move to "In Translation"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("C9").Value = "In Translation"
